$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark issue in row 7 (clone chords above melody...) as Resolved
$ws.Range("C7").Value = "Resolved"

# Mark issue in row 10 (chords should be an octave lower...) as Resolved
$ws.Range("C10").Value = "Resolved"

# Add new row 11 for a newly reported issue
$ws.Range("A11").Value = 42602
$ws.Range("A11").NumberFormat = $ws.Range("A10").NumberFormat

$ws.Range("B11").Value = "last and secondLast classes not playing notes"
$ws.Range("B11").Interior.Color = $ws.Range("B4").Interior.Color

$ws.Range("C11").Value = "Pending"

$ws.Range("B11").Select()
